# Auto-generated Excel COM-interop script
# Updates cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the new symbol-list snapshot, while preserving the
# original text (string) cell type for each updated value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.67"  # was 303.95
$ws.Range("E2").Value = "'5.22%"  # was 5.27%
$ws.Range("D3").Value = "'34.83"  # was 34.85
$ws.Range("E3").Value = "'12.24%"  # was 12.25%
$ws.Range("D4").Value = "'5.169"  # was 5.167
$ws.Range("E4").Value = "'4.39%"  # was 4.07%
$ws.Range("D5").Value = "'0.07800"  # was 0.07822
$ws.Range("E5").Value = "'6.00%"  # was 6.43%
$ws.Range("D6").Value = "'2.298"  # was 2.292
$ws.Range("E6").Value = "'-2.41%"  # was -3.06%
$ws.Range("D7").Value = "'8.051"  # was 8.057
$ws.Range("E7").Value = "'4.23%"  # was 4.24%
$ws.Range("D8").Value = "'3.987"  # was 3.986
$ws.Range("E8").Value = "'6.63%"  # was 6.89%
$ws.Range("D9").Value = "'0.9245"  # was 0.9233
$ws.Range("E9").Value = "'1.33%"  # was 1.26%
$ws.Range("D10").Value = "'0.09987"  # was 0.1004
$ws.Range("E10").Value = "'6.48%"  # was 8.02%
$ws.Range("E11").Value = "'7.99%"  # was 7.83%
$ws.Range("D12").Value = "'0.08543"  # was 0.08618
$ws.Range("E12").Value = "'4.33%"  # was 4.01%
$ws.Range("E13").Value = "'8.18%"  # was 8.20%
$ws.Range("D14").Value = "'0.09915"  # was 0.09911
$ws.Range("E14").Value = "'-0.53%"  # was -0.59%
$ws.Range("D15").Value = "'0.001486"  # was 0.001497
$ws.Range("E15").Value = "'-0.50%"  # was 0.11%
$ws.Range("D16").Value = "'0.04650"  # was 0.04653
$ws.Range("E16").Value = "'2.89%"  # was 2.77%
$ws.Range("D17").Value = "'0.005822"  # was 0.005691
$ws.Range("E17").Value = "'0.76%"  # was -1.10%
$ws.Range("D18").Value = "'3.471"  # was 3.475
$ws.Range("E18").Value = "'0.05%"  # was 0.09%
$ws.Range("D19").Value = "'2.126"  # was 2.127
$ws.Range("E19").Value = "'4.16%"  # was 1.00%
$ws.Range("D20").Value = "'0.3416"  # was 0.3417
$ws.Range("E20").Value = "'2.88%"  # was 2.83%
$ws.Range("D21").Value = "'0.1326"  # was 0.1327
$ws.Range("E21").Value = "'3.03%"  # was 3.10%
$ws.Range("D22").Value = "'4.564"  # was 4.554
$ws.Range("E22").Value = "'9.78%"  # was 9.31%
$ws.Range("D24").Value = "'0.001222"  # was 0.001221
$ws.Range("E24").Value = "'0.70%"  # was 0.75%
$ws.Range("D25").Value = "'0.004333"  # was 0.004332
$ws.Range("E25").Value = "'3.72%"  # was 3.85%
$ws.Range("D26").Value = "'0.0001302"  # was 0.0001301
$ws.Range("E26").Value = "'0.07%"  # was 0.02%
$ws.Range("D27").Value = "'0.0003403"  # was 0.0003400
$ws.Range("E27").Value = "'0.20%"  # was 0.14%
$ws.Range("D39").Value = "'0.01746"  # was 0.01745
$ws.Range("E39").Value = "'11.23%"  # was 10.57%
$ws.Range("E40").Value = "'6.39%"  # was 6.29%
$ws.Range("D41").Value = "'0.007710"  # was 0.007693
$ws.Range("E41").Value = "'4.89%"  # was 4.28%
$ws.Range("D42").Value = "'0.1411"  # was 0.1412
$ws.Range("E42").Value = "'6.06%"  # was 6.11%
$ws.Range("D43").Value = "'0.007650"  # was 0.007261
$ws.Range("E43").Value = "'-22.96%"  # was -26.35%
$ws.Range("D44").Value = "'0.002303"  # was 0.002301
$ws.Range("E44").Value = "'2.30%"  # was 2.70%
$ws.Range("D45").Value = "'0.01002"  # was 0.01001
$ws.Range("E45").Value = "'14.31%"  # was 14.08%
$ws.Range("D46").Value = "'0.00006069"  # was 0.00006066
$ws.Range("E46").Value = "'-0.72%"  # was -0.74%
$ws.Range("D47").Value = "'0.00000000751"  # was 0.00000000750
$ws.Range("E47").Value = "'0.08%"  # was 0.01%
$ws.Range("D48").Value = "'3.880"  # was 5.796
$ws.Range("E48").Value = "'48.87%"  # was 122.41%
$ws.Range("D49").Value = "'0.002694"  # was 0.002691
$ws.Range("E49").Value = "'34.62%"  # was 34.51%
$ws.Range("D50").Value = "'0.00002103"  # was 0.00002101
$ws.Range("E50").Value = "'0.08%"  # was 0.01%
$ws.Range("D51").Value = "'0.0002003"  # was 0.0002001
$ws.Range("E51").Value = "'0.08%"  # was 0.01%
